$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 250.7
$ws.Range("I5").Value = 250.7
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 250.7
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = -135.7
$ws.Range("H6").Value = 151.66667
$ws.Range("I6").Value = 151.66667
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 455.00001
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = -343.00001
$ws.Range("H38").Value = 23
$ws.Range("I38").Value = 23
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 69
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 303
$ws.Range("H43").Value = 8000
$ws.Range("I43").Value = 8000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -7931
$ws.Range("H58").Value = 92.5
$ws.Range("I58").Value = 92.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 277.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -127.5
$ws.Range("H70").Value = 3545.3635
$ws.Range("I70").Value = 3388.889
$ws.Range("J70").Value = 4249.5
$ws.Range("K70").Value = 10166.667
$ws.Range("L70").Value = 12748.5
$ws.Range("M70").Value = -9896.667000000001
$ws.Range("N70").Value = -13288.5
$ws.Range("H73").Value = 3545.3635
$ws.Range("I73").Value = 3388.889
$ws.Range("J73").Value = 4249.5
$ws.Range("K73").Value = 10166.667
$ws.Range("L73").Value = 12748.5
$ws.Range("M73").Value = -9230.667000000001
$ws.Range("N73").Value = -14620.5
$ws.Range("H116").Value = 3495
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 3990
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 3990
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -10874
$ws.Range("H138").Value = 5489.625
$ws.Range("I138").Value = 4237.375
$ws.Range("J138").Value = 5907.0415
$ws.Range("K138").Value = 12712.125
$ws.Range("L138").Value = 17721.1245
$ws.Range("M138").Value = -7572.125
$ws.Range("N138").Value = -28001.1245
$ws.Range("M5").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N38").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2381
$ws.Range("I132").Value = 2381
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7143
$ws.Range("N132").Value = 0
$ws.Range("M132").Value = -4613
$ws.Range("L132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7695.1113
$ws.Range("I31").Value = 4420.6665
$ws.Range("J31").Value = 9332.333000000001
$ws.Range("K31").Value = 4420.6665
$ws.Range("L31").Value = 9332.333000000001
$ws.Range("M31").Value = -4125.6665
$ws.Range("N31").Value = -9922.333000000001
$ws.Range("H34").Value = 7695.1113
$ws.Range("I34").Value = 4420.6665
$ws.Range("J34").Value = 9332.333000000001
$ws.Range("K34").Value = 4420.6665
$ws.Range("L34").Value = 9332.333000000001
$ws.Range("M34").Value = -4218.6665
$ws.Range("N34").Value = -9736.333000000001
$ws.Range("H35").Value = 3059
$ws.Range("I35").Value = 3059
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3059
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = -2765
$ws.Range("H52").Value = 75000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 75000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 75000
$ws.Range("N52").Value = -75588
$ws.Range("H58").Value = 975.5
$ws.Range("I58").Value = 975.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 975.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -772.5
$ws.Range("H105").Value = 3756.8572
$ws.Range("I105").Value = 3574.75
$ws.Range("J105").Value = 3999.6667
$ws.Range("K105").Value = 3574.75
$ws.Range("L105").Value = 3999.6667
$ws.Range("M105").Value = -1827.75
$ws.Range("N105").Value = -7493.6667
$ws.Range("H122").Value = 2999.5
$ws.Range("I122").Value = 2999.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8998.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6548.5
$ws.Range("H136").Value = 975.5
$ws.Range("I136").Value = 975.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2926.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -376.5
$ws.Range("M35").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 172.85715
$ws.Range("I6").Value = 42
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 126
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = -13
$ws.Range("N6").Value = -1726
$ws.Range("H80").Value = 1893.2
$ws.Range("I80").Value = 1903.25
$ws.Range("J80").Value = 1853
$ws.Range("K80").Value = 5709.75
$ws.Range("L80").Value = 5559
$ws.Range("M80").Value = -4773.75
$ws.Range("N80").Value = -7431
$ws.Range("H83").Value = 1893.2
$ws.Range("I83").Value = 1903.25
$ws.Range("J83").Value = 1853
$ws.Range("K83").Value = 17129.25
$ws.Range("L83").Value = 16677
$ws.Range("M83").Value = -12449.25
$ws.Range("N83").Value = -26037
$ws.Range("H122").Value = 987
$ws.Range("I122").Value = 899.3333
$ws.Range("J122").Value = 1250
$ws.Range("K122").Value = 8093.9997
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -5643.9997
$ws.Range("N122").Value = -16150

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H137").Value = 75000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 75000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200
$ws.Range("M132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("H34").Value = 15000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 15000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 15000
$ws.Range("N34").Value = -15344
$ws.Range("H61").Value = 1968.6428
$ws.Range("I61").Value = 1927.7693
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1927.7693
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -1725.7693
$ws.Range("N61").Value = -2904
$ws.Range("H68").Value = 7108.3335
$ws.Range("I68").Value = 3260
$ws.Range("J68").Value = 9857.143
$ws.Range("K68").Value = 3260
$ws.Range("L68").Value = 9857.143
$ws.Range("M68").Value = -2511
$ws.Range("N68").Value = -11355.143
$ws.Range("H71").Value = 7108.3335
$ws.Range("I71").Value = 3260
$ws.Range("J71").Value = 9857.143
$ws.Range("K71").Value = 16300
$ws.Range("L71").Value = 49285.715
$ws.Range("M71").Value = -12556
$ws.Range("N71").Value = -56773.715
$ws.Range("H113").Value = 1968.6428
$ws.Range("I113").Value = 1927.7693
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1927.7693
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 242.2307000000001
$ws.Range("N113").Value = -6840
$ws.Range("L25").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5099.3335
$ws.Range("I62").Value = 4816.5
$ws.Range("J62").Value = 5665
$ws.Range("K62").Value = 4816.5
$ws.Range("L62").Value = 5665
$ws.Range("M62").Value = -4192.5
$ws.Range("N62").Value = -6913
$ws.Range("H65").Value = 5099.3335
$ws.Range("I65").Value = 4816.5
$ws.Range("J65").Value = 5665
$ws.Range("K65").Value = 24082.5
$ws.Range("L65").Value = 28325
$ws.Range("M65").Value = -20962.5
$ws.Range("N65").Value = -34565
